# Applies the "Trade #28 closed at 2026-02-17 20:54:29" update to the
# live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.59   # Current Capital
$summary.Range("B4").Value = 0.38      # Total P&L $
$summary.Range("B5").Value = 0.14      # Total P&L %
$summary.Range("B6").Value = 56        # Total Trades
$summary.Range("B7").Value = 28        # Winning Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking is row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.59     # Capital
$status.Range("D5").Value = 23         # Trades
$status.Range("E5").Value = 0.27       # P&L $
$status.Range("F5").Value = 0.59       # P&L %
$status.Range("G5").Value = 60.87      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - update Trade #28 (row 57) from OPEN to CLOSED
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
# G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
# L Exit Reason, M Duration (min), N Entry Slippage (bps),
# O Exit Slippage (bps), P Confidence, Q Entry Reason
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G57").Value = 0.072049
$allTrades.Range("H57").Value = "CLOSED"
$allTrades.Range("I57").Value = 20.0812
$allTrades.Range("J57").Value = 0.01
$allTrades.Range("K57").Value = 100.59
$allTrades.Range("L57").Value = "early_exit"
$allTrades.Range("M57").Value = 0.15

# New trade row (#89) appended at row 90
$allTrades.Range("A90").Value = 89
# Force text interpretation so the date-like string isn't converted to a
# date serial number, then restore the default "General" number format.
$allTrades.Range("B90").NumberFormat = "@"
$allTrades.Range("B90").Value = "2026-02-17"
$allTrades.Range("B90").NumberFormat = "General"
$allTrades.Range("C90").NumberFormat = "@"
$allTrades.Range("C90").Value = "20:54:23"
$allTrades.Range("C90").NumberFormat = "General"
$allTrades.Range("D90").Value = "MarketMaking"
$allTrades.Range("E90").Value = "UP"
$allTrades.Range("F90").Value = 0.06
$allTrades.Range("H90").Value = "OPEN"
$allTrades.Range("I90").Value = 0
$allTrades.Range("J90").Value = 0
$allTrades.Range("K90").Value = 100.5734535840667
$allTrades.Range("M90").Value = 0
$allTrades.Range("N90").Value = 0
$allTrades.Range("O90").Value = 0
$allTrades.Range("P90").Value = 0.6
$allTrades.Range("Q90").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet - update Trade #28 (row 24) from OPEN to CLOSED
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
# G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
# L Entry Slippage (bps), M Exit Slippage (bps), N Confidence,
# O Entry Reason, P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G24").Value = 0.072049
$marketMaking.Range("H24").Value = "CLOSED"
$marketMaking.Range("I24").Value = 20.0812
$marketMaking.Range("J24").Value = 0.01
$marketMaking.Range("K24").Value = 100.59
$marketMaking.Range("P24").Value = "early_exit"
$marketMaking.Range("Q24").Value = 0.15

# New trade row (#89) appended at row 57
$marketMaking.Range("A57").Value = 89
# Force text interpretation so the date-like string isn't converted to a
# date serial number, then restore the default "General" number format.
$marketMaking.Range("B57").NumberFormat = "@"
$marketMaking.Range("B57").Value = "2026-02-17"
$marketMaking.Range("B57").NumberFormat = "General"
$marketMaking.Range("C57").NumberFormat = "@"
$marketMaking.Range("C57").Value = "20:54:23"
$marketMaking.Range("C57").NumberFormat = "General"
$marketMaking.Range("D57").Value = "MarketMaking"
$marketMaking.Range("E57").Value = "UP"
$marketMaking.Range("F57").Value = 0.06
$marketMaking.Range("H57").Value = "OPEN"
$marketMaking.Range("I57").Value = 0
$marketMaking.Range("J57").Value = 0
$marketMaking.Range("K57").Value = 100.5734535840667
$marketMaking.Range("L57").Value = 0
$marketMaking.Range("M57").Value = 0
$marketMaking.Range("N57").Value = 0.6
$marketMaking.Range("O57").Value = "Normal spread capture: 19600 bps"
$marketMaking.Range("Q57").Value = 0
